# Applies:
#  1) Swap the betting/result data (columns F:V) between row 40 and row 41,
#     leaving columns A:E (Indice, pais, torneio, temporada, data_partida) untouched.
#  2) Append four new match rows (54-57) with data+formats matching the
#     existing rows, extending the sheet from A1:V53 to A1:V57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap row 40 and row 41 (columns F:V only) ---------------------------
$row40 = $ws.Range("F40:V40").Value2
$row41 = $ws.Range("F41:V41").Value2
$ws.Range("F40:V40").Value2 = $row41
$ws.Range("F41:V41").Value2 = $row40

# --- 2) Append new rows 54-57 -----------------------------------------------
# Copy the formatting (styles/number formats) of the last existing row (53)
# onto the new rows, the same way a spreadsheet author extending the table
# downward would.
$ws.Range("A53:V53").Copy() | Out-Null
$ws.Range("A54:V57").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row = 54; A = 53; F = "Floresti";        G = 1; H = "Milsami";          I = 0;
       J = 4.21;  K = "03/12/2023 01:13"; L = 3.84;  M = "03/12/2023 11:41";
       N = 3.68;  O = "03/12/2023 01:13"; P = 3.58;  Q = "03/12/2023 11:41";
       R = 1.65;  S = "03/12/2023 01:13"; T = 1.77;  U = "03/12/2023 11:41";
       V = "https://www.betexplorer.com/football/moldova/super-liga/floresti-milsami/M53yaPf8/" },

    @{ Row = 55; A = 54; F = "Sparta Selemet";  G = 0; H = "Petrocub";         I = 7;
       J = 10.6;  K = "03/12/2023 01:13"; L = 20.49; M = "03/12/2023 11:58";
       N = 6.47;  O = "03/12/2023 01:13"; P = 10.84; Q = "03/12/2023 11:58";
       R = 1.16;  S = "03/12/2023 01:13"; T = 1.03;  U = "03/12/2023 11:58";
       V = "https://www.betexplorer.com/football/moldova/super-liga/sparta-selemet-petrocub-hincesti/tAQlyM21/" },

    @{ Row = 56; A = 55; F = "Balti";           G = 3; H = "Sheriff Tiraspol"; I = 1;
       J = 6.09;  K = "03/12/2023 01:13"; L = 5.71;  M = "03/12/2023 11:55";
       N = 4.12;  O = "03/12/2023 01:13"; P = 3.62;  Q = "03/12/2023 11:55";
       R = 1.42;  S = "03/12/2023 01:13"; T = 1.45;  U = "03/12/2023 11:55";
       V = "https://www.betexplorer.com/football/moldova/super-liga/csf-balti-sheriff-tiraspol/z97ubqAE/" },

    @{ Row = 57; A = 56; F = "Zimbru Chisinau"; G = 1; H = "Dacia Buiucani";   I = 2;
       J = 1.24;  K = "03/12/2023 01:13"; L = 1.24;  M = "03/12/2023 11:46";
       N = 5.22;  O = "03/12/2023 01:13"; P = 4.83;  Q = "03/12/2023 11:46";
       R = 8.33;  S = "03/12/2023 01:13"; T = 7.78;  U = "03/12/2023 11:46";
       V = "https://www.betexplorer.com/football/moldova/super-liga/zimbru-chisinau-dacia-buiucani/K0PhztI7/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A            # Indice
    $ws.Cells.Item($row, 2).Value2 = "moldova"        # pais
    $ws.Cells.Item($row, 3).Value2 = "super-liga"      # torneio
    $ws.Cells.Item($row, 4).Value2 = "2023-2024"       # temporada
    $ws.Cells.Item($row, 5).Value2 = 45263.5           # data_partida
    $ws.Cells.Item($row, 6).Value2 = $r.F               # home
    $ws.Cells.Item($row, 7).Value2 = $r.G               # home_ft_gols
    $ws.Cells.Item($row, 8).Value2 = $r.H               # away
    $ws.Cells.Item($row, 9).Value2 = $r.I               # away_ft_gols
    $ws.Cells.Item($row, 10).Value2 = $r.J              # home_opening_odds
    $ws.Cells.Item($row, 11).Value2 = $r.K              # home_opening_data_hora
    $ws.Cells.Item($row, 12).Value2 = $r.L              # home_closing_odds
    $ws.Cells.Item($row, 13).Value2 = $r.M              # home_closing_data_hora
    $ws.Cells.Item($row, 14).Value2 = $r.N              # draw_opening_odds
    $ws.Cells.Item($row, 15).Value2 = $r.O              # draw_opening_data_hora
    $ws.Cells.Item($row, 16).Value2 = $r.P              # draw_closing_odds
    $ws.Cells.Item($row, 17).Value2 = $r.Q              # draw_closing_data_hora
    $ws.Cells.Item($row, 18).Value2 = $r.R              # away_opening_odds
    $ws.Cells.Item($row, 19).Value2 = $r.S              # away_opening_data_hora
    $ws.Cells.Item($row, 20).Value2 = $r.T              # away_closing_odds
    $ws.Cells.Item($row, 21).Value2 = $r.U              # away_closing_data_hora
    $ws.Cells.Item($row, 22).Value2 = $r.V              # url_partida
}
